$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the new I0 and IF columns
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy formatting (bold font, thin border, centered alignment) from the
# existing header cell H1 onto the two new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for columns I (I0) and J (IF), rows 2-19
$data = @(
    @(2, 1, 3),
    @(3, 1, 5),
    @(4, 1, 6),
    @(5, 1, 5),
    @(6, 1, 4),
    @(7, 1, 5),
    @(8, 1, 5),
    @(9, 1, 6),
    @(10, 4, 7),
    @(11, 1, 4),
    @(12, 1, 6),
    @(13, 1, 6),
    @(14, 4, 7),
    @(15, 2, 7),
    @(16, 5, 5),
    @(17, 1, 3),
    @(18, 1, 3),
    @(19, 1, 2)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($row, 9).Value = $iVal
    $ws.Cells.Item($row, 10).Value = $jVal
}
